$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1538537.5
$ws.Range("I5").Value = 2000096.8
$ws.Range("J5").Value = 7
$ws.Range("K5").Value = 2000096.8
$ws.Range("L5").Value = 7
$ws.Range("M5").Value = -1999981.8
$ws.Range("N5").Value = -237
$ws.Range("H9").Value = 625375.5
$ws.Range("I9").Value = 177.33333
$ws.Range("J9").Value = 1000494.4
$ws.Range("K9").Value = 177.33333
$ws.Range("L9").Value = 1000494.4
$ws.Range("M9").Value = -8.333329999999989
$ws.Range("N9").Value = -1000832.4
$ws.Range("H10").Value = 2158.3333
$ws.Range("I10").Value = 500
$ws.Range("J10").Value = 2490
$ws.Range("K10").Value = 500
$ws.Range("L10").Value = 2490
$ws.Range("M10").Value = -207
$ws.Range("N10").Value = -3076
$ws.Range("H15").Value = 1170.75
$ws.Range("I15").Value = 1170.75
$ws.Range("K15").Value = 3512.25
$ws.Range("M15").Value = -3343.25
$ws.Range("H20").Value = 32500
$ws.Range("I20").Value = 32500
$ws.Range("K20").Value = 32500
$ws.Range("M20").Value = -32270
$ws.Range("H35").Value = 32500
$ws.Range("I35").Value = 32500
$ws.Range("K35").Value = 32500
$ws.Range("M35").Value = -32121
$ws.Range("H41").Value = 979.2
$ws.Range("I41").Value = 862.6
$ws.Range("J41").Value = 1095.8
$ws.Range("K41").Value = 862.6
$ws.Range("L41").Value = 1095.8
$ws.Range("M41").Value = -422.6
$ws.Range("N41").Value = -1975.8
$ws.Range("H43").Value = 1499.2858
$ws.Range("I43").Value = 1415.8334
$ws.Range("K43").Value = 1415.8334
$ws.Range("M43").Value = -1346.8334
$ws.Range("H53").Value = 783.9375
$ws.Range("I53").Value = 644.0769
$ws.Range("J53").Value = 1390
$ws.Range("K53").Value = 644.0769
$ws.Range("L53").Value = 1390
$ws.Range("M53").Value = -7.076900000000023
$ws.Range("N53").Value = -2664
$ws.Range("H96").Value = 938.3333
$ws.Range("I96").Value = 946
$ws.Range("J96").Value = 900
$ws.Range("K96").Value = 2838
$ws.Range("L96").Value = 2700
$ws.Range("M96").Value = -1465
$ws.Range("N96").Value = -5446
$ws.Range("H98").Value = 3878.6
$ws.Range("I98").Value = 2200
$ws.Range("K98").Value = 2200
$ws.Range("M98").Value = -702
$ws.Range("H99").Value = 459
$ws.Range("J99").Value = 670
$ws.Range("L99").Value = 2010
$ws.Range("N99").Value = -5006
$ws.Range("H106").Value = 1844.1428
$ws.Range("I106").Value = 1651.5
$ws.Range("K106").Value = 1651.5
$ws.Range("M106").Value = -1020.5
$ws.Range("H113").Value = 11102.75
$ws.Range("I113").Value = 8268.25
$ws.Range("J113").Value = 13937.25
$ws.Range("K113").Value = 8268.25
$ws.Range("L113").Value = 13937.25
$ws.Range("M113").Value = -5014.25
$ws.Range("N113").Value = -20445.25
$ws.Range("H122").Value = 3878.6
$ws.Range("I122").Value = 2200
$ws.Range("K122").Value = 6600
$ws.Range("M122").Value = -4150
$ws.Range("H135").Value = 951
$ws.Range("I135").Value = 1003.05554
$ws.Range("J135").Value = 763.6
$ws.Range("K135").Value = 9027.49986
$ws.Range("L135").Value = 6872.400000000001
$ws.Range("M135").Value = -6492.49986
$ws.Range("N135").Value = -11942.4
$ws.Range("H137").Value = 2404.75
$ws.Range("I137").Value = 1281.75
$ws.Range("K137").Value = 3845.25
$ws.Range("M137").Value = -1295.25
$ws.Range("H138").Value = 2486.3494
$ws.Range("I138").Value = 4021.0667
$ws.Range("J138").Value = 2147.8088
$ws.Range("K138").Value = 12063.2001
$ws.Range("L138").Value = 6443.426399999999
$ws.Range("M138").Value = -6923.2001
$ws.Range("N138").Value = -16723.4264
$ws.Range("H141").Value = 5838.2856
$ws.Range("I141").Value = 3973.9
$ws.Range("J141").Value = 10499.25
$ws.Range("K141").Value = 11921.7
$ws.Range("L141").Value = 31497.75
$ws.Range("M141").Value = -6741.700000000001
$ws.Range("N141").Value = -41857.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 8181.8184
$ws.Range("I6").Value = 15000
$ws.Range("K6").Value = 15000
$ws.Range("M6").Value = -14827
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H32").Value = 3704.3333
$ws.Range("J32").Value = 9999
$ws.Range("L32").Value = 9999
$ws.Range("N32").Value = -10573
$ws.Range("H61").Value = 2761.8635
$ws.Range("I61").Value = 2844.2
$ws.Range("J61").Value = 2585.4285
$ws.Range("K61").Value = 2844.2
$ws.Range("L61").Value = 2585.4285
$ws.Range("M61").Value = -2632.2
$ws.Range("N61").Value = -3009.4285
$ws.Range("H74").Value = 1789.4762
$ws.Range("I74").Value = 1789.4762
$ws.Range("K74").Value = 1789.4762
$ws.Range("M74").Value = -915.4762000000001
$ws.Range("H77").Value = 1789.4762
$ws.Range("I77").Value = 1789.4762
$ws.Range("K77").Value = 8947.381000000001
$ws.Range("M77").Value = -4579.381000000001
$ws.Range("H97").Value = 1287
$ws.Range("J97").Value = 2912
$ws.Range("L97").Value = 2912
$ws.Range("N97").Value = -3904
$ws.Range("H102").Value = 2159.8
$ws.Range("I102").Value = 2099.75
$ws.Range("K102").Value = 2099.75
$ws.Range("M102").Value = -477.75
$ws.Range("H122").Value = 3531.2222
$ws.Range("I122").Value = 3472.625
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 10417.875
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -7967.875
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 2124.75
$ws.Range("I132").Value = 1239.6
$ws.Range("K132").Value = 3718.8
$ws.Range("M132").Value = -1188.8
$ws.Range("H136").Value = 2761.8635
$ws.Range("I136").Value = 2844.2
$ws.Range("J136").Value = 2585.4285
$ws.Range("K136").Value = 8532.599999999999
$ws.Range("L136").Value = 7756.2855
$ws.Range("M136").Value = -5982.599999999999
$ws.Range("N136").Value = -12856.2855

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4082.1428
$ws.Range("I105").Value = 3595
$ws.Range("K105").Value = 3595
$ws.Range("M105").Value = -1848
$ws.Range("H107").Value = 865.1111
$ws.Range("I107").Value = 865.1111
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 865.1111
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1054.8889
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 2149.3
$ws.Range("I134").Value = 1527.5714
$ws.Range("K134").Value = 4582.7142
$ws.Range("M134").Value = -2047.7142

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H22").Value = 493.33334
$ws.Range("I22").Value = 445
$ws.Range("J22").Value = 590
$ws.Range("K22").Value = 445
$ws.Range("L22").Value = 590
$ws.Range("M22").Value = -95
$ws.Range("N22").Value = -1290
$ws.Range("H31").Value = 2624.2666
$ws.Range("I31").Value = 1583.2222
$ws.Range("J31").Value = 4185.8335
$ws.Range("K31").Value = 1583.2222
$ws.Range("L31").Value = 4185.8335
$ws.Range("M31").Value = -1288.2222
$ws.Range("N31").Value = -4775.8335
$ws.Range("H34").Value = 2624.2666
$ws.Range("I34").Value = 1583.2222
$ws.Range("J34").Value = 4185.8335
$ws.Range("K34").Value = 1583.2222
$ws.Range("L34").Value = 4185.8335
$ws.Range("M34").Value = -1381.2222
$ws.Range("N34").Value = -4589.8335
$ws.Range("H39").Value = 1000
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H49").Value = 1000
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H58").Value = 2274.182
$ws.Range("I58").Value = 2266.8823
$ws.Range("J58").Value = 2299
$ws.Range("K58").Value = 2266.8823
$ws.Range("L58").Value = 2299
$ws.Range("M58").Value = -2063.8823
$ws.Range("N58").Value = -2705
$ws.Range("H88").Value = 17918.428
$ws.Range("J88").Value = 20519.666
$ws.Range("L88").Value = 20519.666
$ws.Range("N88").Value = -21331.666
$ws.Range("H91").Value = 17918.428
$ws.Range("J91").Value = 20519.666
$ws.Range("L91").Value = 20519.666
$ws.Range("N91").Value = -23327.666
$ws.Range("H99").Value = 1499.5
$ws.Range("I99").Value = 1499.5
$ws.Range("K99").Value = 1499.5
$ws.Range("M99").Value = -1.5
$ws.Range("H122").Value = 1657.8334
$ws.Range("J122").Value = 1998
$ws.Range("L122").Value = 5994
$ws.Range("N122").Value = -10894
$ws.Range("H126").Value = 1499.5
$ws.Range("I126").Value = 1499.5
$ws.Range("K126").Value = 4498.5
$ws.Range("M126").Value = -2028.5
$ws.Range("H132").Value = 1693.3448
$ws.Range("I132").Value = 1696.6786
$ws.Range("K132").Value = 5090.0358
$ws.Range("M132").Value = -2560.0358
$ws.Range("H134").Value = 1892.1
$ws.Range("I134").Value = 1938.4375
$ws.Range("J134").Value = 1706.75
$ws.Range("K134").Value = 5815.3125
$ws.Range("L134").Value = 5120.25
$ws.Range("M134").Value = -3280.3125
$ws.Range("N134").Value = -10190.25
$ws.Range("H136").Value = 2274.182
$ws.Range("I136").Value = 2266.8823
$ws.Range("J136").Value = 2299
$ws.Range("K136").Value = 6800.646900000001
$ws.Range("L136").Value = 6897
$ws.Range("M136").Value = -4250.646900000001
$ws.Range("N136").Value = -11997

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2490.3333
$ws.Range("J39").Value = 2490.3333
$ws.Range("L39").Value = 7470.999899999999
$ws.Range("N39").Value = -8058.999899999999
$ws.Range("H44").Value = 1514
$ws.Range("J44").Value = 1376
$ws.Range("L44").Value = 4128
$ws.Range("N44").Value = -4924
$ws.Range("H62").Value = 9599.6
$ws.Range("J62").Value = 2999
$ws.Range("L62").Value = 8997
$ws.Range("N62").Value = -10369
$ws.Range("H65").Value = 9599.6
$ws.Range("J65").Value = 2999
$ws.Range("L65").Value = 26991
$ws.Range("N65").Value = -33855
$ws.Range("H69").Value = 13405.2
$ws.Range("J69").Value = 16009
$ws.Range("L69").Value = 48027
$ws.Range("N69").Value = -49649
$ws.Range("H72").Value = 13405.2
$ws.Range("J72").Value = 16009
$ws.Range("L72").Value = 144081
$ws.Range("N72").Value = -152193
$ws.Range("H122").Value = 477
$ws.Range("I122").Value = 216.5
$ws.Range("J122").Value = 998
$ws.Range("K122").Value = 1948.5
$ws.Range("L122").Value = 8982
$ws.Range("M122").Value = 501.5
$ws.Range("N122").Value = -13882
$ws.Range("H131").Value = 15996.82
$ws.Range("I131").Value = 159667
$ws.Range("J131").Value = 1832.1549
$ws.Range("K131").Value = 479001
$ws.Range("L131").Value = 5496.4647
$ws.Range("M131").Value = -473961
$ws.Range("N131").Value = -15576.4647
$ws.Range("H139").Value = 3729.6538
$ws.Range("I139").Value = 2383.0417
$ws.Range("J139").Value = 19889
$ws.Range("K139").Value = 7149.125100000001
$ws.Range("L139").Value = 59667
$ws.Range("M139").Value = -2009.125100000001
$ws.Range("N139").Value = -69947

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 23290454
$ws.Range("I11").Value = 31860000
$ws.Range("K11").Value = 31860000
$ws.Range("M11").Value = -31859861
$ws.Range("H17").Value = 998
$ws.Range("I17").Value = 998
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 998
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -830
$ws.Range("N17").ClearContents()
$ws.Range("H97").Value = 814.2917
$ws.Range("J97").Value = 661.9091
$ws.Range("L97").Value = 661.9091
$ws.Range("N97").Value = -1653.9091
$ws.Range("H107").Value = 643.0909
$ws.Range("J107").Value = 1025.6364
$ws.Range("L107").Value = 1025.6364
$ws.Range("N107").Value = -4865.6364
$ws.Range("H113").Value = 3545.5833
$ws.Range("I113").Value = 2149.6667
$ws.Range("K113").Value = 2149.6667
$ws.Range("M113").Value = 20.33329999999978
$ws.Range("H126").Value = 3015.353
$ws.Range("I126").Value = 3092.7144
$ws.Range("J126").Value = 2961.2
$ws.Range("K126").Value = 9278.143199999999
$ws.Range("L126").Value = 8883.599999999999
$ws.Range("M126").Value = -6808.143199999999
$ws.Range("N126").Value = -13823.6
$ws.Range("H132").Value = 3050
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2279.0527
$ws.Range("I7").Value = 1478.5714
$ws.Range("J7").Value = 2746
$ws.Range("K7").Value = 1478.5714
$ws.Range("L7").Value = 2746
$ws.Range("M7").Value = -1366.5714
$ws.Range("N7").Value = -2970
$ws.Range("H21").Value = 14338
$ws.Range("J21").Value = 14338
$ws.Range("L21").Value = 14338
$ws.Range("N21").Value = -14686
$ws.Range("H40").Value = 3392.1
$ws.Range("I40").Value = 2998.6667
$ws.Range("J40").Value = 3982.25
$ws.Range("K40").Value = 2998.6667
$ws.Range("L40").Value = 3982.25
$ws.Range("M40").Value = -2862.6667
$ws.Range("N40").Value = -4254.25
$ws.Range("H43").Value = 3264762
$ws.Range("J43").Value = 6196363.5
$ws.Range("L43").Value = 6196363.5
$ws.Range("N43").Value = -6196749.5
$ws.Range("H93").Value = 926.0769
$ws.Range("I93").Value = 922
$ws.Range("K93").Value = 922
$ws.Range("M93").Value = 326
$ws.Range("H122").Value = 3685.25
$ws.Range("I122").Value = 3605.5557
$ws.Range("J122").Value = 3924.3333
$ws.Range("K122").Value = 10816.6671
$ws.Range("L122").Value = 11772.9999
$ws.Range("M122").Value = -8366.667099999999
$ws.Range("N122").Value = -16672.9999
$ws.Range("H126").Value = 2279.0527
$ws.Range("I126").Value = 1478.5714
$ws.Range("J126").Value = 2746
$ws.Range("K126").Value = 4435.7142
$ws.Range("L126").Value = 8238
$ws.Range("M126").Value = -1965.7142
$ws.Range("N126").Value = -13178
$ws.Range("H132").Value = 2282.15
$ws.Range("I132").Value = 2209.7334
$ws.Range("J132").Value = 2499.4
$ws.Range("K132").Value = 6629.2002
$ws.Range("L132").Value = 7498.200000000001
$ws.Range("M132").Value = -4099.2002
$ws.Range("N132").Value = -12558.2
$ws.Range("H136").Value = 1305.8536
$ws.Range("I136").Value = 1203.8649
$ws.Range("J136").Value = 2249.25
$ws.Range("K136").Value = 3611.5947
$ws.Range("L136").Value = 6747.75
$ws.Range("M136").Value = -1061.5947
$ws.Range("N136").Value = -11847.75
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 48999.5
$ws.Range("J15").Value = 48000
$ws.Range("L15").Value = 48000
$ws.Range("N15").Value = -48576
$ws.Range("H43").Value = 28339.25
$ws.Range("J43").Value = 30760
$ws.Range("L43").Value = 30760
$ws.Range("N43").Value = -31058
$ws.Range("H45").Value = 22772
$ws.Range("I45").Value = 22620
$ws.Range("K45").Value = 22620
$ws.Range("M45").Value = -22129
$ws.Range("H74").Value = 23624.666
$ws.Range("J74").Value = 29750.334
$ws.Range("L74").Value = 29750.334
$ws.Range("N74").Value = -31622.334
$ws.Range("H77").Value = 23624.666
$ws.Range("J77").Value = 29750.334
$ws.Range("L77").Value = 89251.00199999999
$ws.Range("N77").Value = -98611.00199999999
$ws.Range("H80").Value = 21665
$ws.Range("J80").Value = 21665
$ws.Range("L80").Value = 21665
$ws.Range("N80").Value = -23661
$ws.Range("H81").Value = 2470.8572
$ws.Range("I81").Value = 2216.1667
$ws.Range("K81").Value = 4432.3334
$ws.Range("M81").Value = -3371.3334
$ws.Range("H83").Value = 21665
$ws.Range("J83").Value = 21665
$ws.Range("L83").Value = 64995
$ws.Range("N83").Value = -74979
$ws.Range("H84").Value = 2470.8572
$ws.Range("I84").Value = 2216.1667
$ws.Range("K84").Value = 22161.667
$ws.Range("M84").Value = -16857.667
$ws.Range("H107").Value = 692.625
$ws.Range("I107").Value = 697.0833
$ws.Range("J107").Value = 688.1667
$ws.Range("K107").Value = 2091.2499
$ws.Range("L107").Value = 2064.5001
$ws.Range("M107").Value = -171.2498999999998
$ws.Range("N107").Value = -5904.5001
$ws.Range("H113").Value = 399.5
$ws.Range("I113").Value = 341
$ws.Range("J113").Value = 545.75
$ws.Range("K113").Value = 1023
$ws.Range("L113").Value = 1637.25
$ws.Range("M113").Value = 1147
$ws.Range("N113").Value = -5977.25
$ws.Range("H122").Value = 6162.391
$ws.Range("I122").Value = 7507.231
$ws.Range("K122").Value = 22521.693
$ws.Range("M122").Value = -20071.693
$ws.Range("H126").Value = 4681.7144
$ws.Range("I126").Value = 4885.091
$ws.Range("J126").Value = 3936
$ws.Range("K126").Value = 14655.273
$ws.Range("L126").Value = 11808
$ws.Range("M126").Value = -12185.273
$ws.Range("N126").Value = -16748
$ws.Range("H132").Value = 1325.9459
$ws.Range("I132").Value = 1278.5518
$ws.Range("K132").Value = 3835.6554
$ws.Range("M132").Value = -1305.6554
$ws.Range("H136").Value = 4348.393
$ws.Range("I136").Value = 3094.0476
$ws.Range("J136").Value = 8111.4287
$ws.Range("K136").Value = 9282.1428
$ws.Range("L136").Value = 24334.2861
$ws.Range("M136").Value = -6732.1428
$ws.Range("N136").Value = -29434.2861
